$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data (columns F:V) between rows 21 and 22.
#    (columns A:E - index/country/tournament/season/date - stay put)
# ---------------------------------------------------------------------------
$ws.Range("F21:V21").Copy()
$ws.Range("F1000:V1000").PasteSpecial(-4163)

$ws.Range("F22:V22").Copy()
$ws.Range("F21:V21").PasteSpecial(-4163)

$ws.Range("F1000:V1000").Copy()
$ws.Range("F22:V22").PasteSpecial(-4163)

$ws.Range("F1000:V1000").Clear()

# ---------------------------------------------------------------------------
# 2) Swap the match data (columns F:V) between rows 81 and 82.
# ---------------------------------------------------------------------------
$ws.Range("F81:V81").Copy()
$ws.Range("F1000:V1000").PasteSpecial(-4163)

$ws.Range("F82:V82").Copy()
$ws.Range("F81:V81").PasteSpecial(-4163)

$ws.Range("F1000:V1000").Copy()
$ws.Range("F82:V82").PasteSpecial(-4163)

$ws.Range("F1000:V1000").Clear()

# ---------------------------------------------------------------------------
# 3) Append three new match rows (89, 90, 91) after the existing last row (88)
#    Start by cloning row 88's formatting down across the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A88:V88").Copy()
$ws.Range("A89:V91").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 89
$ws.Range("A89").Value = 88
$ws.Range("B89").Value = "denmark"
$ws.Range("C89").Value = "superliga"
$ws.Range("D89").Value = "2023-2024"
$ws.Range("E89").Value = 45242.66666666666
$ws.Range("F89").Value = "Midtjylland"
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = "Nordsjaelland"
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2.48
$ws.Range("K89").Value = "05/11/2023 18:13"
$ws.Range("L89").Value = 2.34
$ws.Range("M89").Value = "12/11/2023 15:59"
$ws.Range("N89").Value = 3.46
$ws.Range("O89").Value = "05/11/2023 18:13"
$ws.Range("P89").Value = 3.55
$ws.Range("Q89").Value = "12/11/2023 15:59"
$ws.Range("R89").Value = 2.85
$ws.Range("S89").Value = "05/11/2023 18:13"
$ws.Range("T89").Value = 3.09
$ws.Range("U89").Value = "12/11/2023 15:59"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/denmark/superliga/midtjylland-nordsjaelland/lIzt34bj/"

# Row 90
$ws.Range("A90").Value = 89
$ws.Range("B90").Value = "denmark"
$ws.Range("C90").Value = "superliga"
$ws.Range("D90").Value = "2023-2024"
$ws.Range("E90").Value = 45242.75
$ws.Range("F90").Value = "Odense"
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = "Hvidovre IF"
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 1.56
$ws.Range("K90").Value = "05/11/2023 18:13"
$ws.Range("L90").Value = 1.54
$ws.Range("M90").Value = "12/11/2023 17:52"
$ws.Range("N90").Value = 4.38
$ws.Range("O90").Value = "05/11/2023 18:13"
$ws.Range("P90").Value = 4.49
$ws.Range("Q90").Value = "12/11/2023 17:57"
$ws.Range("R90").Value = 5.05
$ws.Range("S90").Value = "05/11/2023 18:13"
$ws.Range("T90").Value = 6.1
$ws.Range("U90").Value = "12/11/2023 17:57"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/denmark/superliga/odense-hvidovre-if/GOSRQ7jc/"

# Row 91
$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "denmark"
$ws.Range("C91").Value = "superliga"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45242.83333333334
$ws.Range("F91").Value = "Aarhus"
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = "Viborg"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1.67
$ws.Range("K91").Value = "06/11/2023 19:12"
$ws.Range("L91").Value = 2.17
$ws.Range("M91").Value = "12/11/2023 19:58"
$ws.Range("N91").Value = 3.78
$ws.Range("O91").Value = "06/11/2023 19:12"
$ws.Range("P91").Value = 3.28
$ws.Range("Q91").Value = "12/11/2023 19:58"
$ws.Range("R91").Value = 4.96
$ws.Range("S91").Value = "06/11/2023 19:12"
$ws.Range("T91").Value = 3.73
$ws.Range("U91").Value = "12/11/2023 19:58"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/denmark/superliga/aarhus-viborg/6goy4pqp/"
